# The "F1" and "F2" cross labels in column A were ambiguous (each used twice
# for two different rows/environments). Disambiguate them as F1a/F1b and
# F2a/F2b respectively, matching the new source data naming convention.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "F1a"
$ws.Range("A5").Value = "F1b"
$ws.Range("A6").Value = "F2a"
$ws.Range("A7").Value = "F2b"

# Update the sheet's active selection to match the author's saved state.
[void]$ws.Range("E2").Select()
